$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete rows 11-13 (the "Resolving-Mac" sending-cluster rows) entirely.
$ws.Rows.Item(11).Resize(3).EntireRow.Delete()

# Updated numeric values for columns G:T, rows 2-10 (new TPM-based computation).
$values = @{
    2  = @{ G=7.844453333333334;  H=23.53336;          I=0.1489290605659587; J=0.1489290605659588;
             K=2;                 L=0.6666666666666666; M=0.1591403333333333; N=0.477421;
             O=0.01260326986877891; P=0.01260326986877891; Q=1.248368918284444;  R=11.23532026456;
             S=0.001876993141616497; T=0.001876993141616497 }
    3  = @{ G=7.844453333333334;  H=23.53336;          I=0.1489290605659587; J=0.1489290605659588;
             K=3;                 L=1;                  M=3.794584666666667;  N=11.383754;
             O=0.3005157372251983;  P=0.3005157372251983;  Q=29.76644233704889;  R=267.89798103344;
             S=0.04475552643023529; T=0.0447555264302353 }
    4  = @{ G=7.844453333333334;  H=23.53336;          I=0.1489290605659587; J=0.1489290605659588;
             K=3;                 L=1;                  M=8.673183333333334;  N=26.01955;
             O=0.6868809929060228;  P=0.6868809929060229;  Q=68.03638190977779;  R=612.3274371880001;
             S=0.1022965409941069;  T=0.102296540994107 }
    5  = @{ G=20.35396833333334; H=61.06190500000001; I=0.3864255740794268; J=0.3864255740794268;
             K=2;                 L=0.6666666666666666; M=0.1591403333333333; N=0.477421;
             O=0.01260326986877891; P=0.01260326986877891; Q=3.239137305222779;  R=29.152235747005;
             S=0.004870225794320832; T=0.004870225794320832 }
    6  = @{ G=20.35396833333334; H=61.06190500000001; I=0.3864255740794268; J=0.3864255740794268;
             K=3;                 L=1;                  M=3.794584666666667;  N=11.383754;
             O=0.3005157372251983;  P=0.3005157372251983;  Q=77.23485614348557;  R=695.1137052913701;
             S=0.1161269662771494;  T=0.1161269662771494 }
    7  = @{ G=20.35396833333334; H=61.06190500000001; I=0.3864255740794268; J=0.3864255740794268;
             K=3;                 L=1;                  M=8.673183333333334;  N=26.01955;
             O=0.6868809929060228;  P=0.6868809929060229;  Q=176.5336989158612;  R=1588.80329024275;
             S=0.2654283820079565;  T=0.2654283820079566 }
    8  = @{ G=24.47399366666667; H=73.421981;          I=0.4646453653546145; J=0.4646453653546145;
             K=2;                 L=0.6666666666666666; M=0.1591403333333333; N=0.477421;
             O=0.01260326986877891; P=0.01260326986877891; Q=3.894799510111222;  R=35.053195591001;
             S=0.00585605093284158; T=0.005856050932841582 }
    9  = @{ G=24.47399366666667; H=73.421981;          I=0.4646453653546145; J=0.4646453653546145;
             K=3;                 L=1;                  M=3.794584666666667;  N=11.383754;
             O=0.3005157372251983;  P=0.3005157372251983;  Q=92.86864109963045;  R=835.817769896674;
             S=0.1396332445178136;  T=0.1396332445178136 }
    10 = @{ G=24.47399366666667; H=73.421981;          I=0.4646453653546145; J=0.4646453653546145;
             K=3;                 L=1;                  M=8.673183333333334;  N=26.01955;
             O=0.6868809929060228;  P=0.6868809929060229;  Q=212.2674339698389;  R=1910.40690572855;
             S=0.3191560699039593;  T=0.3191560699039594 }
}

$cols = @("G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

foreach ($r in $values.Keys) {
    $rowVals = $values[$r]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value = $rowVals[$c]
    }
}
